# Update the "Förändrad" (Changed) date column (C) from 2023-09-19 (45188)
# to 2023-09-20 (45189) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
